$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 26)
$ws.Range("D2").Value = [double]"0.9999991487265741"
$ws.Range("E2").Value = [double]"0.9999991487265741"

# Row 3 (Control 33)
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = [double]"0.5671078505259338"
$ws.Range("E3").Value = [double]"0.5671078505259338"

# Row 4 (Control 36)
$ws.Range("D4").Value = [double]"0.9999999550198206"
$ws.Range("E4").Value = [double]"0.9999999550198206"

# Row 5 (Control 49)
$ws.Range("D5").Value = [double]"0.0001101138337851883"
$ws.Range("E5").Value = [double]"0.0001101138337851883"

# Row 6 (Control 2)
$ws.Range("D6").Value = [double]"1.136845098984072E-14"
$ws.Range("E6").Value = [double]"1.136845098984072E-14"

# Row 7 (MDD 37)
$ws.Range("D7").Value = [double]"0.9999999999705649"
$ws.Range("E7").Value = [double]"2.943512100728185E-11"

# Row 8 (MDD 24)
$ws.Range("D8").Value = [double]"0.0003717194553826953"
$ws.Range("E8").Value = [double]"0.9996282805446173"

# Row 9 (MDD 6)
$ws.Range("D9").Value = [double]"4.656889161996679E-08"
$ws.Range("E9").Value = [double]"0.9999999534311084"

# Row 10 (MDD 54)
$ws.Range("D10").Value = [double]"1.42208307459904E-12"
$ws.Range("E10").Value = [double]"0.9999999999985779"

# Row 11 (MDD 21)
$ws.Range("D11").Value = [double]"3.069329227767317E-06"
$ws.Range("E11").Value = [double]"0.9999969306707722"
$ws.Range("F11").Value = [double]"9.648359298706055"
$ws.Range("G11").Value = [double]"0.3"
